# Added myself to roster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tim Kukowski's row (row 10) did not yet have a "Subgroup" (column D) value.
# Fill it in with "Communications" (reusing the existing shared string used
# in D3/D11 for the same subgroup).
$ws.Range("D10").Value = "Communications"

# Update the active cell/selection to reflect where editing left off.
$ws.Range("D11").Select()
